# AP_TestData_Login_Nav_Logout.xlsx - "Add files via upload" re-upload.
# The sheet previously held a sample login row (UserName/Password/URL with a
# live hyperlink on the URL cell); the uploaded replacement clears that
# sample data out of row 2, leaving only the header row and the cell
# formatting/style on A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# Remove the hyperlink that lived on A2 (URL sample value).
foreach ($hl in $ws.Hyperlinks) {
    $hl.Delete()
}

# Clear the sample UserName/Password/URL values out of row 2 (A2:C2),
# keeping A2's existing cell style.
$ws.Range("A2:C2").ClearContents()

# Reselect A2:C2 to match the saved selection state.
$ws.Range("A2:C2").Select() | Out-Null
